$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Update the "Hermes Advanced" training entries to "Hermes 5.1 Advanced"
$ws.Range("B15").Value = "Hermes 5.1 Advanced"
$ws.Range("D15").Value = "Projektmanagement Weiterbildung (Hermes 5.1 Advanced)"

# Restore the last active selection recorded in the workbook
$ws.Range("B12").Select()
